$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "79.582.17"
$ws.Range("E2").Value = "  +4.06%  "
$ws.Range("D3").Value = "3.205.07"
$ws.Range("E3").Value = "  +5.32%  "
$ws.Range("D5").Value = "'205.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.39%  "
$ws.Range("D6").Value = "'634.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.27%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  +13.75%  "
$ws.Range("E9").Value = "  +6.02%  "
$ws.Range("D10").Value = "3.209.23"
$ws.Range("E10").Value = "  +5.46%  "
$ws.Range("D11").Value = "'0.592"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +34.74%  "
$ws.Range("E12").Value = "  +3.17%  "
$ws.Range("D13").Value = "'5.51"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.66%  "
$ws.Range("D14").Value = "3.781.59"
$ws.Range("E14").Value = "  +4.84%  "
$ws.Range("E15").Value = "  +18.46%  "
$ws.Range("D17").Value = "79.393.53"
$ws.Range("E17").Value = "  +3.93%  "
$ws.Range("D18").Value = "3.196.13"
$ws.Range("E18").Value = "  +5.24%  "
$ws.Range("D19").Value = "'14.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.57%  "
$ws.Range("D20").Value = "'3.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +31.13%  "
$ws.Range("D21").Value = "'9.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.19%  "
$ws.Range("D22").Value = "'428.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +14.32%  "
$ws.Range("D23").Value = "'5.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +16.79%  "
$ws.Range("D24").Value = "'11.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +13.08%  "
$ws.Range("D25").Value = "3.359.50"
$ws.Range("E25").Value = "  +4.89%  "
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").Value = "'77.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.87%  "
$ws.Range("B27").Value = "NEARProtocol"
$ws.Range("C27").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D27").Value = "'4.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.18%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("E29").Value = "  +7.43%  "
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Value = "'1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.57%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "'9.03"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.88%  "
$ws.Range("D32").Value = "'1.49"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.66%  "
$ws.Range("D33").Value = "'524.33"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.27%  "
$ws.Range("E34").Value = "  +2.23%  "
$ws.Range("E35").Value = "  +27.60%  "
$ws.Range("D36").Value = "'22.90"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.64%  "
$ws.Range("D37").Value = "'0.119"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +12.37%  "
$ws.Range("D38").Value = "'1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.16%  "
$ws.Range("E39").Value = "  +5.35%  "
$ws.Range("D40").Value = "'165.07"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.09%  "
$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D42").Value = "'1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "'192.02"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.56%  "
$ws.Range("E44").Value = "  +6.26%  "
$ws.Range("D45").Value = "'0.819"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -9.95%  "
$ws.Range("E46").Value = "  +7.90%  "
$ws.Range("D47").Value = "'1.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.84%  "
$ws.Range("D48").Value = "'43.10"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.15%  "
$ws.Range("D49").Value = "'25.81"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +15.09%  "
$ws.Range("D50").Value = "'0.639"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.98%  "
$ws.Range("D51").Value = "'2.50"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.40%  "
